# Workbook/sheet handles (provided by the host)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook-level changes -------------------------------------------------
# Remove workbook protection (xl/workbook.xml: <workbookProtection/> dropped)
$wb.Unprotect() | Out-Null

# Rename the (only) worksheet: "Sheet" -> "注文履歴"
$ws.Name = "注文履歴"

# Re-localise the built-in "Normal" cell style's font to the Japanese UI
# default (MS P Gothic) - this is what drives the font-table change in
# xl/styles.xml.
$wb.Styles.Item("Normal").Font.Name = "ＭＳ Ｐゴシック"

# --- row 3: drop the stray empty B3/C3/D3 placeholder cells ----------------
$ws.Range("B3:D3").ClearContents() | Out-Null

# --- new order rows (4-25) --------------------------------------------------
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = '2025-05-15 12:15:19'
$ws.Range("C4").Value = '誠''s唐揚げ'
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1000

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = '2025-05-15 12:15:19'
$ws.Range("C5").Value = 'test'
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 111

$ws.Range("A6").Value = '合計'
$ws.Range("E6").Value = 1111

$ws.Range("A7").Value = '合計'
$ws.Range("E7").Value = 0

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = '2025-05-16 01:13:42'
$ws.Range("C8").Value = '誠''s唐揚げ'
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1000

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = '2025-05-16 01:13:42'
$ws.Range("C9").Value = 'test'
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 111

$ws.Range("A10").Value = '合計'
$ws.Range("E10").Value = 1111

$ws.Range("A11").Value = 7
$ws.Range("B11").Value = '2025-05-16 03:27:14'
$ws.Range("C11").Value = 'test'
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 111

$ws.Range("A12").Value = '合計'
$ws.Range("E12").Value = 111

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = '2025-05-19 03:00:02'
$ws.Range("C13").Value = '誠''s唐揚げ'
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 4000

$ws.Range("A14").Value = '合計'
$ws.Range("E14").Value = 4000

$ws.Range("A15").Value = '完了'
$ws.Range("F15").Value = '2025-05-21 06:02:17'

$ws.Range("A16").Value = '完了'
$ws.Range("F16").Value = '2025-05-21 06:02:20'

$ws.Range("A17").Value = '完了'
$ws.Range("F17").Value = '2025-05-21 06:02:21'

$ws.Range("A18").Value = '完了'
$ws.Range("F18").Value = '2025-05-21 06:02:22'

$ws.Range("A19").Value = '完了'
$ws.Range("F19").Value = '2025-05-21 06:02:24'

$ws.Range("A20").Value = '完了'
$ws.Range("F20").Value = '2025-05-21 06:02:25'

$ws.Range("A21").Value = '完了'
$ws.Range("F21").Value = '2025-05-21 06:02:26'

$ws.Range("A22").Value = '完了'
$ws.Range("F22").Value = '2025-05-21 06:02:30'

$ws.Range("A23").Value = 9
$ws.Range("B23").Value = '2025-05-22 09:03:56'
$ws.Range("C23").Value = '誠''s唐揚げ'
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 2000

$ws.Range("A24").Value = 9
$ws.Range("B24").Value = '2025-05-22 09:03:56'
$ws.Range("C24").Value = 'test'
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 111

$ws.Range("A25").Value = '合計'
$ws.Range("E25").Value = 2111

# --- row 26: final "完了" marker row (B26:E26 are intentionally left blank,
# matching the source which carries no real value in those columns)
$ws.Range("A26").Value = '完了'
$ws.Range("F26").Value = '2025-05-22 09:04:08'

# --- sheet view: select F27 (one past the last data row) and make this the
# active/visible tab -----------------------------------------------------
$ws.Range("F27").Select() | Out-Null
